$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H69").Value = 3818.3928
$ws.Range("I69").Value = 3880.6667
$ws.Range("J69").Value = 3788.8948
$ws.Range("K69").Value = 11642.0001
$ws.Range("L69").Value = 11366.6844
$ws.Range("M69").Value = -10768.0001
$ws.Range("N69").Value = -13114.6844
$ws.Range("H72").Value = 3818.3928
$ws.Range("I72").Value = 3880.6667
$ws.Range("J72").Value = 3788.8948
$ws.Range("K72").Value = 34926.0003
$ws.Range("L72").Value = 34100.0532
$ws.Range("M72").Value = -30558.0003
$ws.Range("N72").Value = -42836.0532
$ws.Range("H101").Value = 400
$ws.Range("I101").Value = 400
$ws.Range("K101").Value = 1200
$ws.Range("M101").Value = 422
$ws.Range("H112").Value = 6050.516
$ws.Range("J112").Value = 7426.64
$ws.Range("L112").Value = 22279.92
$ws.Range("N112").Value = -24495.92
$ws.Range("H113").Value = 2944.3333
$ws.Range("I113").Value = 2760
$ws.Range("J113").Value = 3174.75
$ws.Range("K113").Value = 2760
$ws.Range("L113").Value = 3174.75
$ws.Range("M113").Value = 494
$ws.Range("N113").Value = -9682.75
$ws.Range("H129").Value = 1013.6842
$ws.Range("J129").Value = 1083.25
$ws.Range("L129").Value = 3249.75
$ws.Range("N129").Value = -13249.75
$ws.Range("H138").Value = 2119.1804
$ws.Range("I138").Value = 1474.4412
$ws.Range("J138").Value = 2931.074
$ws.Range("K138").Value = 4423.3236
$ws.Range("L138").Value = 8793.222
$ws.Range("M138").Value = 716.6764000000003
$ws.Range("N138").Value = -19073.222

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 547102.9399999999
$ws.Range("I32").Value = 641746
$ws.Range("J32").Value = 17101.8
$ws.Range("K32").Value = 641746
$ws.Range("L32").Value = 17101.8
$ws.Range("M32").Value = -641459
$ws.Range("N32").Value = -17675.8
$ws.Range("H61").Value = 1752.381
$ws.Range("I61").Value = 1387.4166
$ws.Range("J61").Value = 2920.2666
$ws.Range("K61").Value = 1387.4166
$ws.Range("L61").Value = 2920.2666
$ws.Range("M61").Value = -1175.4166
$ws.Range("N61").Value = -3344.2666
$ws.Range("H74").Value = 1188.6
$ws.Range("I74").Value = 869
$ws.Range("J74").Value = 1801.1666
$ws.Range("K74").Value = 869
$ws.Range("L74").Value = 1801.1666
$ws.Range("M74").Value = 5
$ws.Range("N74").Value = -3549.1666
$ws.Range("H77").Value = 1188.6
$ws.Range("I77").Value = 869
$ws.Range("J77").Value = 1801.1666
$ws.Range("K77").Value = 4345
$ws.Range("L77").Value = 9005.833000000001
$ws.Range("M77").Value = 23
$ws.Range("N77").Value = -17741.833
$ws.Range("H110").Value = 1113.8334
$ws.Range("I110").Value = 1099.9231
$ws.Range("J110").Value = 1150
$ws.Range("K110").Value = 1099.9231
$ws.Range("L110").Value = 1150
$ws.Range("M110").Value = 945.0769
$ws.Range("N110").Value = -5240
$ws.Range("H136").Value = 1752.381
$ws.Range("I136").Value = 1387.4166
$ws.Range("J136").Value = 2920.2666
$ws.Range("K136").Value = 4162.2498
$ws.Range("L136").Value = 8760.799800000001
$ws.Range("M136").Value = -1612.2498
$ws.Range("N136").Value = -13860.7998

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 736.3333
$ws.Range("I16").Value = 467.33334
$ws.Range("K16").Value = 467.33334
$ws.Range("M16").Value = -180.33334
$ws.Range("H31").Value = 4882.73
$ws.Range("I31").Value = 1257.4166
$ws.Range("J31").Value = 9716.481
$ws.Range("K31").Value = 1257.4166
$ws.Range("L31").Value = 9716.481
$ws.Range("M31").Value = -962.4166
$ws.Range("N31").Value = -10306.481
$ws.Range("H34").Value = 4882.73
$ws.Range("I34").Value = 1257.4166
$ws.Range("J34").Value = 9716.481
$ws.Range("K34").Value = 1257.4166
$ws.Range("L34").Value = 9716.481
$ws.Range("M34").Value = -1055.4166
$ws.Range("N34").Value = -10120.481
$ws.Range("H113").Value = 736.3333
$ws.Range("I113").Value = 467.33334
$ws.Range("K113").Value = 467.33334
$ws.Range("M113").Value = 1702.66666

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 806.5789
$ws.Range("I5").Value = 570.8333
$ws.Range("J5").Value = 1210.7142
$ws.Range("K5").Value = 1712.4999
$ws.Range("L5").Value = 3632.1426
$ws.Range("M5").Value = -1600.4999
$ws.Range("N5").Value = -3856.1426
$ws.Range("H117").Value = 532.8570999999999
$ws.Range("I117").Value = 293.33334
$ws.Range("J117").Value = 712.5
$ws.Range("K117").Value = 880.0000200000001
$ws.Range("L117").Value = 2137.5
$ws.Range("M117").Value = 2561.99998
$ws.Range("N117").Value = -9021.5
$ws.Range("H131").Value = 1330.3
$ws.Range("J131").Value = 1706.6923
$ws.Range("L131").Value = 5120.0769
$ws.Range("N131").Value = -15200.0769
$ws.Range("H132").Value = 2318.7021
$ws.Range("J132").Value = 2317.361
$ws.Range("L132").Value = 20856.249
$ws.Range("N132").Value = -25916.249
$ws.Range("H135").Value = 806.5789
$ws.Range("I135").Value = 570.8333
$ws.Range("J135").Value = 1210.7142
$ws.Range("K135").Value = 5137.4997
$ws.Range("L135").Value = 10896.4278
$ws.Range("M135").Value = -2602.4997
$ws.Range("N135").Value = -15966.4278
$ws.Range("H137").Value = 6674644.5
$ws.Range("I137").Value = 23824466
$ws.Range("J137").Value = 5269.8335
$ws.Range("K137").Value = 71473398
$ws.Range("L137").Value = 15809.5005
$ws.Range("M137").Value = -71468298
$ws.Range("N137").Value = -26009.5005
$ws.Range("H140").Value = 1858.8235
$ws.Range("I140").Value = 1580
$ws.Range("J140").Value = 2172.5
$ws.Range("K140").Value = 4740
$ws.Range("L140").Value = 6517.5
$ws.Range("M140").Value = 440
$ws.Range("N140").Value = -16877.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 101919.3
$ws.Range("I113").Value = 201478.6
$ws.Range("J113").Value = 2360
$ws.Range("K113").Value = 201478.6
$ws.Range("L113").Value = 2360
$ws.Range("M113").Value = -199308.6
$ws.Range("N113").Value = -6700
$ws.Range("H132").Value = 1844.7733
$ws.Range("I132").Value = 1529.4286
$ws.Range("J132").Value = 3500.3333
$ws.Range("K132").Value = 4588.2858
$ws.Range("L132").Value = 10500.9999
$ws.Range("M132").Value = -2058.2858
$ws.Range("N132").Value = -15560.9999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 17859030
$ws.Range("I16").Value = 2000.3334
$ws.Range("J16").Value = 28573246
$ws.Range("K16").Value = 2000.3334
$ws.Range("L16").Value = 28573246
$ws.Range("M16").Value = -1830.3334
$ws.Range("N16").Value = -28573586
$ws.Range("H61").Value = 5071.4287
$ws.Range("I61").Value = 5933.3335
$ws.Range("J61").Value = 4425
$ws.Range("K61").Value = 5933.3335
$ws.Range("L61").Value = 4425
$ws.Range("M61").Value = -5731.3335
$ws.Range("N61").Value = -4829
$ws.Range("H113").Value = 5071.4287
$ws.Range("I113").Value = 5933.3335
$ws.Range("J113").Value = 4425
$ws.Range("K113").Value = 5933.3335
$ws.Range("L113").Value = 4425
$ws.Range("M113").Value = -3763.3335
$ws.Range("N113").Value = -8765

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 5563.222
$ws.Range("I81").Value = 6011.6665
$ws.Range("J81").Value = 4666.3335
$ws.Range("K81").Value = 12023.333
$ws.Range("L81").Value = 9332.666999999999
$ws.Range("M81").Value = -10962.333
$ws.Range("N81").Value = -11454.667
$ws.Range("H84").Value = 5563.222
$ws.Range("I84").Value = 6011.6665
$ws.Range("J84").Value = 4666.3335
$ws.Range("K84").Value = 60116.665
$ws.Range("L84").Value = 46663.335
$ws.Range("M84").Value = -54812.665
$ws.Range("N84").Value = -57271.335
$ws.Range("H122").Value = 2214.68
$ws.Range("I122").Value = 2197.5293
$ws.Range("J122").Value = 2251.125
$ws.Range("K122").Value = 6592.5879
$ws.Range("L122").Value = 6753.375
$ws.Range("M122").Value = -4142.5879
$ws.Range("N122").Value = -11653.375
